$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 335, shifting existing rows 335:418 down to 336:419
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new weekly data record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R repeat the same boilerplate values used
# throughout this sheet (Perejil / Vega Central Mapocho de Santiago).
$ws.Range("A335").Value2 = 9
$ws.Range("B335").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C335").Value = "Metropolitana"
$ws.Range("D335").Value2 = 44782
$ws.Range("E335").Value2 = 13
$ws.Range("F335").Value2 = 100112044
$ws.Range("G335").Value = "Perejil"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value2 = 82
$ws.Range("K335").Value2 = 18000
$ws.Range("L335").Value2 = 19000
$ws.Range("M335").Value2 = 18390
$ws.Range("N335").Value = "`$/docena de atados"
$ws.Range("O335").Value = "Región Metropolitana"
$ws.Range("P335").Value2 = 6130
$ws.Range("Q335").Value2 = 3
$ws.Range("R335").Value = "Hortaliza"
